# "Finish Story folder structure, Finish Ogier story file"
#
# The source workbook is being promoted from a plain story-text sheet to
# the finished "Deepdive" file, so the sheet is renamed, and the author's
# leftover working position (scrolled down to row ~395 with G6 selected)
# is cleared back to the top of the sheet (A1) before handing the file off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the only worksheet: "Ogier Story" -> "Ogier Deepdive"
$ws.Name = "Ogier Deepdive"

# Make sure the sheet is the active one, then reset the view/selection back
# to the top-left corner (A1) instead of the mid-document working position
# (topLeftCell="A395", selection G6) that was left over from editing.
$ws.Activate()
[void]$ws.Range("A1").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
